$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# API routes changed for Excel Sheet.
# Update Endpoint (column C) values to reflect new route naming.
$ws.Range("C8").Value  = "decksvote/:id/upvote"
$ws.Range("C9").Value  = "decksvote/:id/downvote"
$ws.Range("C10").Value = "flashcardsvote/:id/upvote"
$ws.Range("C11").Value = "flashcardsvote/:id/downvote"

$ws.Range("C16").Value = "flashcards"
$ws.Range("C17").Value = "flashcards/:id"
$ws.Range("C18").Value = "flashcards"
$ws.Range("C19").Value = "flashcards/:id"
$ws.Range("C20").Value = "flashcards/:id"

$ws.Range("C21").Value = "decks"
$ws.Range("C22").Value = "decks/:id"
$ws.Range("C23").Value = "decks"
$ws.Range("C24").Value = "decks/:id"
$ws.Range("C25").Value = "decks/:id"

# Update the saved view state (active selection).
$ws.Range("C16").Select()
